# Update the cached "datetimeFigureOut" field text (Insert > Header & Footer
# "Date and time" placeholder) from 7/6/2018 to 4/2/19 everywhere it is
# defined: the slide master, every slide layout, and the notes master.

$p = $ppt.ActivePresentation
$newDate = "4/2/19"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# Slide master.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every slide layout hanging off the slide master.
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# Notes master.
$notesMaster = $p.NotesMaster
Update-DatePlaceholder $notesMaster.Shapes
